$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the time-range text in B5 to add the second session and wrap the text
$ws.Range("B5").Value = "9.00-10.45, 11:45-13.00"
$ws.Range("B5").WrapText = $true

# Let Excel auto-fit row 5's height to the newly wrapped two-line text
$ws.Rows.Item(5).AutoFit()

# Update the active selection to match the authored state
$ws.Range("B5").Select()
